$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 6333.3335
$ws.Cells.Item(58, 9).Value = 2000
$ws.Cells.Item(58, 10).Value = 7200
$ws.Cells.Item(58, 11).Value = 6000
$ws.Cells.Item(58, 12).Value = 21600
$ws.Cells.Item(58, 13).Value = -5850
$ws.Cells.Item(58, 14).Value = -21900

$ws.Cells.Item(69, 8).Value = 5000
$ws.Cells.Item(69, 10).Value = 5000
$ws.Cells.Item(69, 12).Value = 15000
$ws.Cells.Item(69, 14).Value = -16748

$ws.Cells.Item(72, 8).Value = 5000
$ws.Cells.Item(72, 10).Value = 5000
$ws.Cells.Item(72, 12).Value = 45000
$ws.Cells.Item(72, 14).Value = -53736

$ws.Cells.Item(98, 8).Value = 962.5833
$ws.Cells.Item(98, 9).Value = 1006
$ws.Cells.Item(98, 11).Value = 1006
$ws.Cells.Item(98, 13).Value = 492

$ws.Cells.Item(122, 8).Value = 962.5833
$ws.Cells.Item(122, 9).Value = 1006
$ws.Cells.Item(122, 11).Value = 3018
$ws.Cells.Item(122, 13).Value = -568

$ws.Cells.Item(132, 8).Value = 4980.6665
$ws.Cells.Item(132, 9).Value = 2902
$ws.Cells.Item(132, 11).Value = 8706
$ws.Cells.Item(132, 13).Value = -6176

$ws.Cells.Item(137, 8).Value = 3347.3872
$ws.Cells.Item(137, 9).Value = 2620.0908
$ws.Cells.Item(137, 10).Value = 3747.4
$ws.Cells.Item(137, 11).Value = 7860.2724
$ws.Cells.Item(137, 12).Value = 11242.2
$ws.Cells.Item(137, 13).Value = -5310.2724
$ws.Cells.Item(137, 14).Value = -16342.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2831.0386
$ws.Cells.Item(45, 9).Value = 2708.48
$ws.Cells.Item(45, 11).Value = 2708.48
$ws.Cells.Item(45, 13).Value = -2331.48

$ws.Cells.Item(61, 8).Value = 3749.5
$ws.Cells.Item(61, 9).Value = 3749.5
$ws.Cells.Item(61, 11).Value = 3749.5
$ws.Cells.Item(61, 13).Value = -3537.5

$ws.Cells.Item(74, 8).Value = 2268.4285
$ws.Cells.Item(74, 9).Value = 1979.8334
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 11).Value = 1979.8334
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 13).Value = -1105.8334
$ws.Cells.Item(74, 14).Value = -5748

$ws.Cells.Item(77, 8).Value = 2268.4285
$ws.Cells.Item(77, 9).Value = 1979.8334
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 11).Value = 9899.166999999999
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = -5531.166999999999
$ws.Cells.Item(77, 14).Value = -28736

$ws.Cells.Item(110, 8).Value = 1215.2
$ws.Cells.Item(110, 9).Value = 1215.2
$ws.Cells.Item(110, 11).Value = 1215.2
$ws.Cells.Item(110, 13).Value = 829.8

$ws.Cells.Item(122, 8).Value = 1611.6666
$ws.Cells.Item(122, 9).Value = 1611.6666
$ws.Cells.Item(122, 11).Value = 4834.9998
$ws.Cells.Item(122, 13).Value = -2384.9998

$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 3749.5
$ws.Cells.Item(136, 9).Value = 3749.5
$ws.Cells.Item(136, 11).Value = 11248.5
$ws.Cells.Item(136, 13).Value = -8698.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2612.5
$ws.Cells.Item(105, 9).Value = 2380
$ws.Cells.Item(105, 11).Value = 2380
$ws.Cells.Item(105, 13).Value = -633

$ws.Cells.Item(134, 8).Value = 3251.125
$ws.Cells.Item(134, 9).Value = 3251.125
$ws.Cells.Item(134, 11).Value = 9753.375
$ws.Cells.Item(134, 13).Value = -7218.375

$ws.Cells.Item(138, 8).Value = 50000
$ws.Cells.Item(138, 10).Value = 50000
$ws.Cells.Item(138, 12).Value = 50000
$ws.Cells.Item(138, 14).Value = -60280

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 619.8
$ws.Cells.Item(16, 9).Value = 524.75
$ws.Cells.Item(16, 11).Value = 524.75
$ws.Cells.Item(16, 13).Value = -237.75

$ws.Cells.Item(39, 8).Value = 17499.5
$ws.Cells.Item(39, 9).Value = 5000
$ws.Cells.Item(39, 10).Value = 29999
$ws.Cells.Item(39, 11).Value = 5000
$ws.Cells.Item(39, 12).Value = 29999
$ws.Cells.Item(39, 13).Value = -4609
$ws.Cells.Item(39, 14).Value = -30781

$ws.Cells.Item(49, 8).Value = 17499.5
$ws.Cells.Item(49, 9).Value = 5000
$ws.Cells.Item(49, 10).Value = 29999
$ws.Cells.Item(49, 11).Value = 5000
$ws.Cells.Item(49, 12).Value = 29999
$ws.Cells.Item(49, 13).Value = -4818
$ws.Cells.Item(49, 14).Value = -30363

$ws.Cells.Item(58, 8).Value = 3495.6667
$ws.Cells.Item(58, 9).Value = 3495.6667
$ws.Cells.Item(58, 11).Value = 3495.6667
$ws.Cells.Item(58, 13).Value = -3292.6667

$ws.Cells.Item(105, 8).Value = 1237.5714
$ws.Cells.Item(105, 9).Value = 1237.5714
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 1237.5714
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = 509.4286
$ws.Cells.Item(105, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 619.8
$ws.Cells.Item(113, 9).Value = 524.75
$ws.Cells.Item(113, 11).Value = 524.75
$ws.Cells.Item(113, 13).Value = 1645.25

$ws.Cells.Item(122, 8).Value = 1998.8
$ws.Cells.Item(122, 9).Value = 1998.8
$ws.Cells.Item(122, 11).Value = 5996.4
$ws.Cells.Item(122, 13).Value = -3546.4

$ws.Cells.Item(132, 8).Value = 2239.1333
$ws.Cells.Item(132, 9).Value = 2257
$ws.Cells.Item(132, 11).Value = 6771
$ws.Cells.Item(132, 13).Value = -4241

$ws.Cells.Item(136, 8).Value = 3495.6667
$ws.Cells.Item(136, 9).Value = 3495.6667
$ws.Cells.Item(136, 11).Value = 10487.0001
$ws.Cells.Item(136, 13).Value = -7937.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 1809.75
$ws.Cells.Item(109, 9).Value = 1809.75
$ws.Cells.Item(109, 11).Value = 5429.25
$ws.Cells.Item(109, 13).Value = -4389.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()

$ws.Cells.Item(102, 8).Value = 5750
$ws.Cells.Item(102, 9).Value = 5750
$ws.Cells.Item(102, 11).Value = 5750
$ws.Cells.Item(102, 13).Value = -4128

$ws.Cells.Item(123, 8).Value = 60000
$ws.Cells.Item(123, 10).Value = 60000
$ws.Cells.Item(123, 12).Value = 60000
$ws.Cells.Item(123, 14).Value = -64900

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4242.4443
$ws.Cells.Item(7, 10).Value = 4518.2
$ws.Cells.Item(7, 12).Value = 4518.2
$ws.Cells.Item(7, 14).Value = -4742.2

$ws.Cells.Item(46, 8).Value = 891.5
$ws.Cells.Item(46, 10).Value = 937.25
$ws.Cells.Item(46, 12).Value = 937.25
$ws.Cells.Item(46, 14).Value = -1313.25

$ws.Cells.Item(55, 8).Value = 6187.125
$ws.Cells.Item(55, 9).Value = 5999.4287
$ws.Cells.Item(55, 11).Value = 5999.4287
$ws.Cells.Item(55, 13).Value = -5826.4287

$ws.Cells.Item(68, 8).Value = 2963.2727
$ws.Cells.Item(68, 9).Value = 2844.111
$ws.Cells.Item(68, 11).Value = 2844.111
$ws.Cells.Item(68, 13).Value = -2095.111

$ws.Cells.Item(71, 8).Value = 2963.2727
$ws.Cells.Item(71, 9).Value = 2844.111
$ws.Cells.Item(71, 11).Value = 14220.555
$ws.Cells.Item(71, 13).Value = -10476.555

$ws.Cells.Item(122, 8).Value = 5750
$ws.Cells.Item(122, 9).Value = 4000
$ws.Cells.Item(122, 10).Value = 6333.3335
$ws.Cells.Item(122, 11).Value = 12000
$ws.Cells.Item(122, 12).Value = 19000.0005
$ws.Cells.Item(122, 13).Value = -9550
$ws.Cells.Item(122, 14).Value = -23900.0005

$ws.Cells.Item(126, 8).Value = 4242.4443
$ws.Cells.Item(126, 10).Value = 4518.2
$ws.Cells.Item(126, 12).Value = 13554.6
$ws.Cells.Item(126, 14).Value = -18494.6

$ws.Cells.Item(132, 8).Value = 23686.7
$ws.Cells.Item(132, 9).Value = 23114.625
$ws.Cells.Item(132, 11).Value = 69343.875
$ws.Cells.Item(132, 13).Value = -66813.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 54666
$ws.Cells.Item(2, 9).Value = 54666
$ws.Cells.Item(2, 11).Value = 54666
$ws.Cells.Item(2, 13).Value = -54554

$ws.Cells.Item(96, 8).Value = 750
$ws.Cells.Item(96, 9).Value = 750
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 750
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = 623
$ws.Cells.Item(96, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 510.45456
$ws.Cells.Item(113, 9).Value = 497.22223
$ws.Cells.Item(113, 11).Value = 1491.66669
$ws.Cells.Item(113, 13).Value = 678.33331

$ws.Cells.Item(122, 8).Value = 669166.3
$ws.Cells.Item(122, 9).Value = 2000000
$ws.Cells.Item(122, 11).Value = 6000000
$ws.Cells.Item(122, 13).Value = -5997550

$ws.Cells.Item(136, 8).Value = 1927.7391
$ws.Cells.Item(136, 9).Value = 1746.95
$ws.Cells.Item(136, 10).Value = 3133
$ws.Cells.Item(136, 11).Value = 5240.85
$ws.Cells.Item(136, 12).Value = 9399
$ws.Cells.Item(136, 13).Value = -2690.85
$ws.Cells.Item(136, 14).Value = -14499
